$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 25.38222733333333
$ws.Range("H2").Value = 76.146682
$ws.Range("I2").Value = 0.1760862452187379
$ws.Range("J2").Value = 0.1760862452187379
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 147.4213356666667
$ws.Range("N2").Value = 442.264007
$ws.Range("O2").Value = 0.9507885170992249
$ws.Range("P2").Value = 0.950788517099225
$ws.Range("Q2").Value = 3741.881855674975
$ws.Range("R2").Value = 33676.93670107477
$ws.Range("S2").Value = 0.1674207799730943
$ws.Range("T2").Value = 0.1674207799730943

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 25.38222733333333
$ws.Range("H3").Value = 76.146682
$ws.Range("I3").Value = 0.1760862452187379
$ws.Range("J3").Value = 0.1760862452187379
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.340788333333334
$ws.Range("N3").Value = 7.022365000000001
$ws.Range("O3").Value = 0.01509682881537204
$ws.Range("P3").Value = 0.01509682881537204
$ws.Range("Q3").Value = 59.41442161588112
$ws.Range("R3").Value = 534.72979454293
$ws.Range("S3").Value = 0.002658343900808909
$ws.Range("T3").Value = 0.002658343900808909

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 25.38222733333333
$ws.Range("H4").Value = 76.146682
$ws.Range("I4").Value = 0.1760862452187379
$ws.Range("J4").Value = 0.1760862452187379
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.289533666666667
$ws.Range("N4").Value = 15.868601
$ws.Range("O4").Value = 0.03411465408540306
$ws.Range("P4").Value = 0.03411465408540307
$ws.Range("Q4").Value = 134.2601460146536
$ws.Range("R4").Value = 1208.341314131882
$ws.Range("S4").Value = 0.006007121344834703
$ws.Range("T4").Value = 0.006007121344834705

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 94.773687
$ws.Range("H5").Value = 284.321061
$ws.Range("I5").Value = 0.6574814128880592
$ws.Range("J5").Value = 0.6574814128880593
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 147.4213356666667
$ws.Range("N5").Value = 442.264007
$ws.Range("O5").Value = 0.9507885170992249
$ws.Range("P5").Value = 0.950788517099225
$ws.Range("Q5").Value = 13971.6635235946
$ws.Range("R5").Value = 125744.9717123514
$ws.Range("S5").Value = 0.625125777580141
$ws.Range("T5").Value = 0.6251257775801411

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 94.773687
$ws.Range("H6").Value = 284.321061
$ws.Range("I6").Value = 0.6574814128880592
$ws.Range("J6").Value = 0.6574814128880593
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.340788333333334
$ws.Range("N6").Value = 7.022365000000001
$ws.Range("O6").Value = 0.01509682881537204
$ws.Range("P6").Value = 0.01509682881537204
$ws.Range("Q6").Value = 221.845140836585
$ws.Range("R6").Value = 1996.606267529265
$ws.Range("S6").Value = 0.00992588433965997
$ws.Range("T6").Value = 0.009925884339659971

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 94.773687
$ws.Range("H7").Value = 284.321061
$ws.Range("I7").Value = 0.6574814128880592
$ws.Range("J7").Value = 0.6574814128880593
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 5.289533666666667
$ws.Range("N7").Value = 15.868601
$ws.Range("O7").Value = 0.03411465408540306
$ws.Range("P7").Value = 0.03411465408540307
$ws.Range("Q7").Value = 501.308608100629
$ws.Range("R7").Value = 4511.777472905661
$ws.Range("S7").Value = 0.0224297509682582
$ws.Range("T7").Value = 0.02242975096825821

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 23.99065033333333
$ws.Range("H8").Value = 71.971951
$ws.Range("I8").Value = 0.1664323418932028
$ws.Range("J8").Value = 0.1664323418932028
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 147.4213356666667
$ws.Range("N8").Value = 442.264007
$ws.Range("O8").Value = 0.9507885170992249
$ws.Range("P8").Value = 0.950788517099225
$ws.Range("Q8").Value = 3536.733715651962
$ws.Range("R8").Value = 31830.60344086766
$ws.Range("S8").Value = 0.1582419595459895
$ws.Range("T8").Value = 0.1582419595459895

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 23.99065033333333
$ws.Range("H9").Value = 71.971951
$ws.Range("I9").Value = 0.1664323418932028
$ws.Range("J9").Value = 0.1664323418932028
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.340788333333334
$ws.Range("N9").Value = 7.022365000000001
$ws.Range("O9").Value = 0.01509682881537204
$ws.Range("P9").Value = 0.01509682881537204
$ws.Range("Q9").Value = 56.15703440934612
$ws.Range("R9").Value = 505.4133096841151
$ws.Range("S9").Value = 0.002512600574903154
$ws.Range("T9").Value = 0.002512600574903155

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 23.99065033333333
$ws.Range("H10").Value = 71.971951
$ws.Range("I10").Value = 0.1664323418932028
$ws.Range("J10").Value = 0.1664323418932028
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 5.289533666666667
$ws.Range("N10").Value = 15.868601
$ws.Range("O10").Value = 0.03411465408540306
$ws.Range("P10").Value = 0.03411465408540307
$ws.Range("Q10").Value = 126.8993526233946
$ws.Range("R10").Value = 1142.094173610551
$ws.Range("S10").Value = 0.00567778177231015
$ws.Range("T10").Value = 0.005677781772310152

